# Refresh the cryptocurrency price/volume snapshot on the "cryptos" sheet
# (mirrors the scheduled GitHub Actions scrape update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to text
# so Excel does not silently convert them (losing formatting like trailing zeros
# or the thousand-dot grouping used by this sheet).
$textCells = @("D4", "D5", "D6", "D8", "D10", "D11", "D13", "D14", "D15", "D16", "D18", "D21", "D23", "D25", "D27", "D29", "D30", "D31", "D32", "D33", "D34", "D36", "D37", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D47", "D48", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values scraped for this run
$ws.Range("D2").Value = "29.392.61"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "1.846.66"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "239.92"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("D6").Value = "0.6295"
$ws.Range("E6").Value = "  -1.06%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "0.07589"
$ws.Range("E8").Value = "  +0.39%  "
$ws.Range("E9").Value = "  -1.35%  "
$ws.Range("D10").Value = "24.49"
$ws.Range("E10").Value = "  -0.74%  "
$ws.Range("D11").Value = "0.07738"
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("D12").Value = "1.840.18"
$ws.Range("E12").Value = "  -0.34%  "
$ws.Range("D13").Value = "0.00001092"
$ws.Range("E13").Value = "  +9.49%  "
$ws.Range("D14").Value = "4.998"
$ws.Range("E14").Value = "  -0.03%  "
$ws.Range("D15").Value = "0.6778"
$ws.Range("E15").Value = "  -1.00%  "
$ws.Range("D16").Value = "83.64"
$ws.Range("E16").Value = "  +0.73%  "
$ws.Range("D17").Value = "2.091.43"
$ws.Range("E17").Value = "  -7.64%  "
$ws.Range("D18").Value = "6.145"
$ws.Range("E18").Value = "  -0.58%  "
$ws.Range("D19").Value = "29.410.24"
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("E20").Value = "  -0.52%  "
$ws.Range("D21").Value = "12.44"
$ws.Range("E21").Value = "  -0.34%  "
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").Value = "7.421"
$ws.Range("E23").Value = "  -1.98%  "
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").Value = "157.39"
$ws.Range("E25").Value = "  +0.22%  "
$ws.Range("E26").Value = "  -0.78%  "
$ws.Range("D27").Value = "8.378"
$ws.Range("E27").Value = "  -0.17%  "
$ws.Range("E28").Value = "  -0.29%  "
$ws.Range("D29").Value = "1.463"
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").Value = "1.307"
$ws.Range("E30").Value = "  +4.43%  "
$ws.Range("D31").Value = "0.05602"
$ws.Range("E31").Value = "  -1.96%  "
$ws.Range("D32").Value = "4.102"
$ws.Range("E32").Value = "  -0.60%  "
$ws.Range("D33").Value = "4.033"
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("D34").Value = "1.845"
$ws.Range("E34").Value = "  -0.39%  "
$ws.Range("D36").Value = "0.7098"
$ws.Range("E36").Value = "  -0.92%  "
$ws.Range("D37").Value = "2.583"
$ws.Range("E37").Value = "  -0.43%  "
$ws.Range("D38").Value = "1.229.83"
$ws.Range("E38").Value = "  -1.92%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.01798"
$ws.Range("E39").Value = "  -0.44%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "2.769"
$ws.Range("E40").Value = "  -0.70%  "
$ws.Range("D41").Value = "6.441"
$ws.Range("E41").Value = "  +4.03%  "
$ws.Range("D42").Value = "0.9043"
$ws.Range("E42").Value = "  -0.45%  "
$ws.Range("D43").Value = "1.000"
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("D44").Value = "101.70"
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("D45").Value = "65.94"
$ws.Range("E45").Value = "  -1.00%  "
$ws.Range("E46").Value = "  +3.59%  "
$ws.Range("D47").Value = "7.189"
$ws.Range("E47").Value = "  +1.43%  "
$ws.Range("D48").Value = "0.4018"
$ws.Range("E48").Value = "  -0.27%  "
$ws.Range("D49").Value = "8.990"
$ws.Range("E49").Value = "  -2.25%  "
$ws.Range("E50").Value = "  -1.59%  "
$ws.Range("D51").Value = "0.1120"
$ws.Range("E51").Value = "  -0.88%  "
